# Auto-generated edit script: apply the scraped schedule update to all three sheets.
$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 18:50:27'
$ws.Cells.Item(3, 1).Value = 'Total filas: 290'
$ws.Cells.Item(15, 3).Value = '225_GOMEZ'
$ws.Cells.Item(16, 3).Value = '215A_EL PATO'
$ws.Cells.Item(43, 1).Value = '06:58:58'
$ws.Cells.Item(43, 3).Value = '15_ABASTO'
$ws.Cells.Item(43, 4).Value = 90
$ws.Cells.Item(44, 1).Value = '08:14:55'
$ws.Cells.Item(44, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(44, 4).Value = 14
$ws.Cells.Item(45, 1).Value = '07:51:40'
$ws.Cells.Item(45, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(45, 4).Value = 38
$ws.Cells.Item(46, 1).Value = '08:14:55'
$ws.Cells.Item(46, 3).Value = '15_ABASTO'
$ws.Cells.Item(46, 4).Value = 15
$ws.Cells.Item(71, 1).Value = '08:14:55'
$ws.Cells.Item(71, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(71, 4).Value = 77
$ws.Cells.Item(72, 1).Value = '08:49:06'
$ws.Cells.Item(72, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(72, 4).Value = 42
$ws.Cells.Item(103, 1).Value = '11:01:19'
$ws.Cells.Item(103, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(103, 4).Value = 14
$ws.Cells.Item(104, 1).Value = '10:32:07'
$ws.Cells.Item(104, 3).Value = '14_ABASTO'
$ws.Cells.Item(104, 4).Value = 43
$ws.Cells.Item(117, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(118, 3).Value = '17_ROMERO'
$ws.Cells.Item(148, 1).Value = '12:43:13'
$ws.Cells.Item(148, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(148, 4).Value = 42
$ws.Cells.Item(149, 1).Value = '12:58:23'
$ws.Cells.Item(149, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(149, 4).Value = 27
$ws.Cells.Item(162, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(163, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(199, 1).Value = '16:02:30'
$ws.Cells.Item(199, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(199, 4).Value = 3
$ws.Cells.Item(200, 1).Value = '14:58:38'
$ws.Cells.Item(200, 3).Value = '14_ABASTO'
$ws.Cells.Item(200, 4).Value = 67
$ws.Cells.Item(223, 1).Value = '16:34:05'
$ws.Cells.Item(223, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(223, 4).Value = 33
$ws.Cells.Item(224, 1).Value = '16:57:38'
$ws.Cells.Item(224, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(224, 4).Value = 10
$ws.Cells.Item(230, 1).Value = '16:34:05'
$ws.Cells.Item(230, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(230, 4).Value = 60
$ws.Cells.Item(231, 1).Value = '17:22:11'
$ws.Cells.Item(231, 3).Value = '17_ROMERO'
$ws.Cells.Item(231, 4).Value = 12
$ws.Cells.Item(244, 1).Value = '17:22:11'
$ws.Cells.Item(244, 3).Value = '10_OLMOS'
$ws.Cells.Item(244, 4).Value = 31
$ws.Cells.Item(245, 1).Value = '16:34:05'
$ws.Cells.Item(245, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(245, 4).Value = 79
$ws.Cells.Item(265, 1).Value = '18:50:27'
$ws.Cells.Item(265, 4).Value = 2
$ws.Cells.Item(266, 1).Value = '18:50:27'
$ws.Cells.Item(266, 2).Value = '18:53'
$ws.Cells.Item(266, 3).Value = '14_ABASTO'
$ws.Cells.Item(266, 4).Value = 3
$ws.Cells.Item(267, 1).Value = '18:50:27'
$ws.Cells.Item(267, 2).Value = '18:57'
$ws.Cells.Item(267, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(267, 4).Value = 7
$ws.Cells.Item(268, 1).Value = '18:34:20'
$ws.Cells.Item(268, 2).Value = '18:59'
$ws.Cells.Item(268, 4).Value = 25
$ws.Cells.Item(269, 1).Value = '18:01:05'
$ws.Cells.Item(269, 2).Value = '19:00'
$ws.Cells.Item(269, 3).Value = '14_ABASTO'
$ws.Cells.Item(269, 4).Value = 59
$ws.Cells.Item(270, 1).Value = '18:50:27'
$ws.Cells.Item(270, 2).Value = '19:04'
$ws.Cells.Item(270, 3).Value = '215_EL PELIGRO'
$ws.Cells.Item(270, 4).Value = 14
$ws.Cells.Item(271, 1).Value = '18:50:27'
$ws.Cells.Item(271, 2).Value = '19:10'
$ws.Cells.Item(271, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(271, 4).Value = 20
$ws.Cells.Item(272, 2).Value = '19:16'
$ws.Cells.Item(272, 3).Value = '17_ROMERO'
$ws.Cells.Item(272, 4).Value = 42
$ws.Cells.Item(273, 1).Value = '18:50:27'
$ws.Cells.Item(273, 2).Value = '19:16'
$ws.Cells.Item(273, 4).Value = 26
$ws.Cells.Item(274, 1).Value = '18:50:27'
$ws.Cells.Item(274, 2).Value = '19:17'
$ws.Cells.Item(274, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(274, 4).Value = 27
$ws.Cells.Item(275, 1).Value = '18:34:20'
$ws.Cells.Item(275, 2).Value = '19:17'
$ws.Cells.Item(275, 4).Value = 43
$ws.Cells.Item(276, 1).Value = '18:34:20'
$ws.Cells.Item(276, 2).Value = '19:18'
$ws.Cells.Item(276, 4).Value = 44
$ws.Cells.Item(277, 1).Value = '17:47:31'
$ws.Cells.Item(277, 2).Value = '19:18'
$ws.Cells.Item(277, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(277, 4).Value = 91
$ws.Cells.Item(278, 1).Value = '18:50:27'
$ws.Cells.Item(278, 2).Value = '19:21'
$ws.Cells.Item(278, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(278, 4).Value = 31
$ws.Cells.Item(279, 1).Value = '18:50:27'
$ws.Cells.Item(279, 2).Value = '19:23'
$ws.Cells.Item(279, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(279, 4).Value = 33
$ws.Cells.Item(280, 1).Value = '18:01:05'
$ws.Cells.Item(280, 2).Value = '19:23'
$ws.Cells.Item(280, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(280, 4).Value = 82
$ws.Cells.Item(281, 1).Value = '18:50:27'
$ws.Cells.Item(281, 2).Value = '19:28'
$ws.Cells.Item(281, 3).Value = '215C_EL PATO'
$ws.Cells.Item(281, 4).Value = 38
$ws.Cells.Item(282, 1).Value = '18:50:27'
$ws.Cells.Item(282, 2).Value = '19:35'
$ws.Cells.Item(282, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(282, 4).Value = 45
$ws.Cells.Item(283, 1).Value = '18:01:05'
$ws.Cells.Item(283, 2).Value = '19:36'
$ws.Cells.Item(283, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(283, 4).Value = 95
$ws.Cells.Item(284, 1).Value = '18:50:27'
$ws.Cells.Item(284, 2).Value = '19:39'
$ws.Cells.Item(284, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(284, 4).Value = 49
$ws.Cells.Item(285, 1).Value = '17:47:31'
$ws.Cells.Item(285, 2).Value = '19:40'
$ws.Cells.Item(285, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(285, 4).Value = 113
$ws.Cells.Item(286, 1).Value = '18:50:27'
$ws.Cells.Item(286, 2).Value = '19:52'
$ws.Cells.Item(286, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(286, 4).Value = 62
$ws.Cells.Item(287, 1).Value = '18:50:27'
$ws.Cells.Item(287, 2).Value = '19:53'
$ws.Cells.Item(287, 3).Value = '225_GOMEZ'
$ws.Cells.Item(287, 4).Value = 63
$ws.Cells.Item(288, 1).Value = '18:50:27'
$ws.Cells.Item(288, 2).Value = '20:05'
$ws.Cells.Item(288, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(288, 4).Value = 75
$ws.Cells.Item(288, 5).Value = 'LP1912'
$ws.Cells.Item(289, 1).Value = '18:50:27'
$ws.Cells.Item(289, 2).Value = '20:07'
$ws.Cells.Item(289, 3).Value = '215C_EL PATO'
$ws.Cells.Item(289, 4).Value = 77
$ws.Cells.Item(289, 5).Value = 'LP1912'
$ws.Cells.Item(290, 1).Value = '18:50:27'
$ws.Cells.Item(290, 2).Value = '20:12'
$ws.Cells.Item(290, 3).Value = '14_ABASTO'
$ws.Cells.Item(290, 4).Value = 82
$ws.Cells.Item(290, 5).Value = 'LP1912'
$ws.Cells.Item(291, 1).Value = '18:50:27'
$ws.Cells.Item(291, 2).Value = '20:22'
$ws.Cells.Item(291, 3).Value = '15_ABASTO'
$ws.Cells.Item(291, 4).Value = 92
$ws.Cells.Item(291, 5).Value = 'LP1912'
$ws.Cells.Item(292, 1).Value = '18:50:27'
$ws.Cells.Item(292, 2).Value = '20:30'
$ws.Cells.Item(292, 3).Value = '10_OLMOS'
$ws.Cells.Item(292, 4).Value = 100
$ws.Cells.Item(292, 5).Value = 'LP1912'
$ws.Cells.Item(293, 1).Value = '18:34:20'
$ws.Cells.Item(293, 2).Value = '20:31'
$ws.Cells.Item(293, 3).Value = '10_OLMOS'
$ws.Cells.Item(293, 4).Value = 117
$ws.Cells.Item(293, 5).Value = 'LP1912'
$ws.Cells.Item(294, 1).Value = '18:50:27'
$ws.Cells.Item(294, 2).Value = '20:46'
$ws.Cells.Item(294, 3).Value = '17_ROMERO'
$ws.Cells.Item(294, 4).Value = 116
$ws.Cells.Item(294, 5).Value = 'LP1912'
$ws.Cells.Item(295, 1).Value = '18:50:27'
$ws.Cells.Item(295, 2).Value = '20:47'
$ws.Cells.Item(295, 3).Value = '215B_EL PATO'
$ws.Cells.Item(295, 4).Value = 117
$ws.Cells.Item(295, 5).Value = 'LP1912'

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 18:50:27'
$ws.Cells.Item(3, 1).Value = 'Total filas: 47'
$ws.Cells.Item(49, 1).Value = '18:50:27'
$ws.Cells.Item(49, 4).Value = 14
$ws.Cells.Item(50, 1).Value = '18:50:27'
$ws.Cells.Item(50, 4).Value = 38
$ws.Cells.Item(51, 1).Value = '18:50:27'
$ws.Cells.Item(51, 4).Value = 77
$ws.Cells.Item(52, 1).Value = '18:50:27'
$ws.Cells.Item(52, 2).Value = '20:47'
$ws.Cells.Item(52, 3).Value = '215B_EL PATO'
$ws.Cells.Item(52, 4).Value = 117
$ws.Cells.Item(52, 5).Value = 'LP1912'

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 18:50:27'
$ws.Cells.Item(44, 1).Value = '18:50:27'
$ws.Cells.Item(44, 4).Value = 34
$ws.Cells.Item(45, 1).Value = '18:50:27'
$ws.Cells.Item(45, 4).Value = 68

